$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 4
    5  = 7
    6  = 6
    7  = 5
    8  = 1
    9  = 5
    10 = 4
    11 = 2
    12 = 3
    13 = 3
    14 = 2
    15 = 4
    16 = 2
    17 = 2
    18 = 3
    19 = 6
    20 = 6
    21 = 1
    22 = 4
    23 = 3
    24 = 1
    25 = 5
    26 = 4
    27 = 3
    28 = 5
    29 = 4
    30 = 3
    31 = 2
    32 = 4
    33 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
